{"js": "// Update weekly excess mortality analyses \u2014 replace the text of 8\n// \"SourceCode\"-styled paragraphs (R console output quoted strings) with\n// their updated-week versions, per the commit's unified diff.\n\nconst replacements = [\n  {\n    find: \"CBS heeft het aantal overlijdensgevallen bijgewerkt\",\n    oldText:\n      '## [1] \"CBS heeft het aantal overlijdensgevallen bijgewerkt t/m week 27 van dit jaar. Deze week combineer ik de grafieken over sterfte per week (alleen totaal) met een uitgebreide ondersterfte vergelijking tussen 2018 en 2020. Ik kijk ook naar de oversterfte in Europa.\"',\n    newText:\n      '## [1] \"CBS heeft het aantal overlijdensgevallen bijgewerkt t/m week 28 van dit jaar. Deze week combineer ik de grafieken over sterfte per week (alleen totaal) met een uitgebreide ondersterfte vergelijking tussen 2018 en 2020. Ik kijk ook naar de oversterfte in Europa.\"',\n  },\n  {\n    find: \"Sterfte per week: De blauwe piek\",\n    oldText:\n      '## [1] \"Sterfte per week: De blauwe piek die je ziet is 2020. Gemiddeld aantal overledenen in week 27 (2015-2019) is 2791, 2020 = 2607. RIVM zegt nu 12 in week 27. Er is dus ondersterfte in week 27, zelfs met 12 offici\u00eble corona-overledenen (wat waarschijnlijk niet eens alles is).\"',\n    newText:\n      '## [1] \"Sterfte per week: De blauwe piek die je ziet is 2020. Gemiddeld aantal overledenen in week 28 (2015-2019) is 2684, 2020 = 2581. RIVM zegt nu 5 in week 28. Er is dus ondersterfte in week 28, zelfs met 5 offici\u00eble corona-overledenen (wat waarschijnlijk niet eens alles is).\"',\n  },\n  {\n    find: \"Oversterfte NL: Voor week 12 t/m 19\",\n    oldText:\n      \"## [1] \\\"Oversterfte NL: Voor week 12 t/m 19 voorspelden de verschillende methodes dit: (1) 9236, (2) 8189, en (3) 8617. Ik houd de oversterfte in de 'heftige' periode dus op 8189-9236.\\\"\",\n    newText:\n      \"## [1] \\\"Oversterfte NL: Voor week 12 t/m 19 voorspelden de verschillende methodes dit: (1) 9239, (2) 8192, en (3) 8617. Ik houd de oversterfte in de 'heftige' periode dus op 8192-9239.\\\"\",\n  },\n  {\n    find: \"Op basis van methode (2) is in week\",\n    oldText:\n      '## [1] \"Op basis van methode (2) is in week 27 de ondersterfte -308. Ondersterfte vanaf week 20 t/m 27: (1) -305, (2) -1271 en (3) -1577. Ik houd de totale sterfte (week 12 t/m 27) op dit moment tussen de 6918-8931.\"',\n    newText:\n      '## [1] \"Op basis van methode (2) is in week 28 de ondersterfte -220. Ondersterfte vanaf week 20 t/m 28: (1) -377, (2) -1462 en (3) -1824. Ik houd de totale sterfte (week 12 t/m 28) op dit moment tussen de 6730-8862.\"',\n  },\n  {\n    find: \"De ondersterfte begon in 2018 vanaf week 14\",\n    oldText:\n      '## [1] \"De ondersterfte begon in 2018 vanaf week 14 voor alle leeftijdsgroepen. De ondersterfte van week 14 t/m 27 was -1966 (flink hoger dan in 2020). Maar dat is uiteraard ondersterfte over meer weken, namelijk 13. Dus wat als we net als in 2020 de eerste 8 weken van 2018 pakken?\"',\n    newText:\n      '## [1] \"De ondersterfte begon in 2018 vanaf week 14 voor alle leeftijdsgroepen. De ondersterfte van week 14 t/m 28 was -2003 (flink hoger dan in 2020). Maar dat is uiteraard ondersterfte over meer weken, namelijk 14. Dus wat als we net als in 2020 de eerste 9 weken van 2018 pakken?\"',\n  },\n  {\n    find: \"Als we het aantal weken ondersterfte gelijkzetten\",\n    oldText:\n      \"## [1] \\\"Als we het aantal weken ondersterfte gelijkzetten (zodat we de 'snelheid' waarmee oversterfte 'gecorrigeerd' wordt kunnen vergelijken), gebruiken we nu 8 weken. De ondersterfte in 2020 gaat nu dus sneller dan in 2018 (-1271 vs. -1037). Het is afwachten of deze trend doorzet.\\\"\",\n    newText:\n      \"## [1] \\\"Als we het aantal weken ondersterfte gelijkzetten (zodat we de 'snelheid' waarmee oversterfte 'gecorrigeerd' wordt kunnen vergelijken), gebruiken we nu 9 weken. De ondersterfte in 2020 gaat nu dus sneller dan in 2018 (-1462 vs. -1037). Het is afwachten of deze trend doorzet.\\\"\",\n  },\n  {\n    find: \"Europa week 28: EuroMOMO\",\n    oldText:\n      '## [1] \"Europa week 28: EuroMOMO lijkt wat correcties doorgevoerd te hebben en heeft de oversterfte afgelopen week onderschat. De dalende trend was dus iets te voorbarig, het blijft vrij stabiel. Totale oversterfte week 12 t/m 28 = 169000\"',\n    newText:\n      '## [1] \"Europa week 29: Een aantal landen laten weer wat oversterfte zien (Belgie, Zweden, Spanje, en zeker Portugal), maar Europees breed blijft het gelijk. Totale oversterfte week 12 t/m 29 = 170000\"',\n  },\n  {\n    find: \"Conclusie na week 27\",\n    oldText:\n      '## [1] \"Conclusie na week 27: we zitten in een periode van ondersterfte (dat is normaal na een heftige epidemie, zoals ik deze week cijfermatig laat zien), waarbij de ondersterfte snel toeneemt. Afhankelijk van het type zomer kan dit ook wel flink doorzetten, maar dat is afwachten.\"',\n    newText:\n      '## [1] \"Conclusie na week 28: we zitten in een periode van ondersterfte (dat is normaal na een heftige epidemie, zoals ik deze week cijfermatig laat zien), waarbij de ondersterfte snel toeneemt. Afhankelijk van het type zomer kan dit ook wel flink doorzetten, maar dat is afwachten.\"',\n  },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nfor (const para of paragraphs.items) {\n  for (const rep of replacements) {\n    if (para.text.includes(rep.find)) {\n      para.insertText(rep.newText, \"Replace\");\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update weekly excess mortality analyses \u2014 replace the text of 8\n# \"SourceCode\"-styled paragraphs (R console output quoted strings) with\n# their updated-week versions, per the commit's unified diff.\n#\n# Each paragraph is located by a unique substring of its current text,\n# then the paragraph's Range.Text is assigned the full new text (the\n# paragraph-mark character at Range.End is excluded so the paragraph\n# itself, and its run formatting / style, is preserved).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Find = \"CBS heeft het aantal overlijdensgevallen bijgewerkt\"\n        New  = '## [1] \"CBS heeft het aantal overlijdensgevallen bijgewerkt t/m week 28 van dit jaar. Deze week combineer ik de grafieken over sterfte per week (alleen totaal) met een uitgebreide ondersterfte vergelijking tussen 2018 en 2020. Ik kijk ook naar de oversterfte in Europa.\"'\n    },\n    @{\n        Find = \"Sterfte per week: De blauwe piek\"\n        New  = '## [1] \"Sterfte per week: De blauwe piek die je ziet is 2020. Gemiddeld aantal overledenen in week 28 (2015-2019) is 2684, 2020 = 2581. RIVM zegt nu 5 in week 28. Er is dus ondersterfte in week 28, zelfs met 5 offici\u00eble corona-overledenen (wat waarschijnlijk niet eens alles is).\"'\n    },\n    @{\n        Find = \"Oversterfte NL: Voor week 12 t/m 19\"\n        New  = \"## [1] \"\"Oversterfte NL: Voor week 12 t/m 19 voorspelden de verschillende methodes dit: (1) 9239, (2) 8192, en (3) 8617. Ik houd de oversterfte in de 'heftige' periode dus op 8192-9239.\"\"\"\n    },\n    @{\n        Find = \"Op basis van methode (2) is in week\"\n        New  = '## [1] \"Op basis van methode (2) is in week 28 de ondersterfte -220. Ondersterfte vanaf week 20 t/m 28: (1) -377, (2) -1462 en (3) -1824. Ik houd de totale sterfte (week 12 t/m 28) op dit moment tussen de 6730-8862.\"'\n    },\n    @{\n        Find = \"De ondersterfte begon in 2018 vanaf week 14\"\n        New  = '## [1] \"De ondersterfte begon in 2018 vanaf week 14 voor alle leeftijdsgroepen. De ondersterfte van week 14 t/m 28 was -2003 (flink hoger dan in 2020). Maar dat is uiteraard ondersterfte over meer weken, namelijk 14. Dus wat als we net als in 2020 de eerste 9 weken van 2018 pakken?\"'\n    },\n    @{\n        Find = \"Als we het aantal weken ondersterfte gelijkzetten\"\n        New  = \"## [1] \"\"Als we het aantal weken ondersterfte gelijkzetten (zodat we de 'snelheid' waarmee oversterfte 'gecorrigeerd' wordt kunnen vergelijken), gebruiken we nu 9 weken. De ondersterfte in 2020 gaat nu dus sneller dan in 2018 (-1462 vs. -1037). Het is afwachten of deze trend doorzet.\"\"\"\n    },\n    @{\n        Find = \"Europa week 28: EuroMOMO\"\n        New  = '## [1] \"Europa week 29: Een aantal landen laten weer wat oversterfte zien (Belgie, Zweden, Spanje, en zeker Portugal), maar Europees breed blijft het gelijk. Totale oversterfte week 12 t/m 29 = 170000\"'\n    },\n    @{\n        Find = \"Conclusie na week 27\"\n        New  = '## [1] \"Conclusie na week 28: we zitten in een periode van ondersterfte (dat is normaal na een heftige epidemie, zoals ik deze week cijfermatig laat zien), waarbij de ondersterfte snel toeneemt. Afhankelijk van het type zomer kan dit ook wel flink doorzetten, maar dat is afwachten.\"'\n    }\n)\n\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    $t = $p.Range.Text\n    foreach ($rep in $replacements) {\n        if ($t.Contains($rep.Find)) {\n            $r = $p.Range\n            $r.End = $r.End - 1\n            $r.Text = $rep.New\n            break\n        }\n    }\n}\n"}
